# Insert a new data row at row 689 (pushing the existing rows 689..760 down
# to 690..761) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 689..760 down to 690..761, leaving a blank row 689 behind.
$ws.Rows.Item(689).Insert()

# Populate the newly inserted row 689 with the new record's data.
$ws.Cells.Item(689, 1).Value2 = 4
$ws.Cells.Item(689, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(689, 3).Value = "Los Lagos"
$ws.Cells.Item(689, 4).Value2 = 45194
$ws.Cells.Item(689, 5).Value2 = 10
$ws.Cells.Item(689, 6).Value2 = 100114001
$ws.Cells.Item(689, 7).Value = "Papa"
$ws.Cells.Item(689, 8).Value = "Asterix"
$ws.Cells.Item(689, 9).Value = "1a (guarda)"
$ws.Cells.Item(689, 10).Value2 = 150
$ws.Cells.Item(689, 11).Value2 = 28000
$ws.Cells.Item(689, 12).Value2 = 28000
$ws.Cells.Item(689, 13).Value2 = 28000
$ws.Cells.Item(689, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(689, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(689, 16).Value2 = 1120
$ws.Cells.Item(689, 17).Value2 = 25
$ws.Cells.Item(689, 18).Value = "Hortaliza"
